$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C17 value from 123456 to 12345
$ws.Range("C17").Value = 12345

# Clear the formatting on B17 (removes the extra cellXf with applyFill="1")
$ws.Range("B17").ClearFormats()

# Update the active selection to E13
$ws.Range("E13").Select()
